# Convert some images to code
#
# Remove the shapes that were replaced by the generated-code text boxes:
#   - Freeform 25           (id 26) - background split rectangle behind the two "code" boxes
#   - Content Placeholder 3 (id 34) - right-hand "code" text box
#   - Content Placeholder 3 (id 35) - left-hand "code" text box
#   - Rectangle 35          (id 36) - "(c)" label
#   - Rectangle 36          (id 37) - "(d)" label
#
# Shape id 34 also has an entrance animation (appear) registered in the
# slide's timeline, so remove that effect (which also cleans up the
# corresponding <p:bldP> build entry) before deleting the shape itself.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$idsToRemove = @(26, 34, 35, 36, 37)

# Remove any timeline effects that target one of the shapes we're about to
# delete (iterate backwards since removing shifts indices).
$seq = $s.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $eff = $seq.Item($i)
    if ($idsToRemove -contains $eff.Shape.Id) {
        $eff.Delete()
    }
}

# Now remove the shapes themselves.
foreach ($targetId in $idsToRemove) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shape = $s.Shapes.Item($i)
        if ($shape.Id -eq $targetId) {
            $shape.Delete()
            break
        }
    }
}
